# Update the "want to go" counts (column F) across sheets to reflect the
# latest scraped data (gh-pages output generated at 456a3b4).
#
# NOTE: this engine's named-parameter binding (e.g. "-Foo bar") does not
# reliably pass argument values into functions, so positional parameters
# are used for the helper function below.

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param([string]$SheetName, [hashtable]$RowValues)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F" + $row).Value = $RowValues[$row]
    }
}

# Sheet "展览" (Exhibitions)
Set-FValues "展览" @{
    2  = 7854
    3  = 107
    5  = 13666
    7  = 593
    8  = 648
    9  = 450
    10 = 140
    11 = 441
    12 = 781
    13 = 44
    14 = 85
    15 = 343
    16 = 20
    17 = 298
    18 = 141
    19 = 398
    20 = 192
    21 = 1096
    22 = 80
    23 = 641
    24 = 2233
    25 = 754
    26 = 55
    27 = 555
    28 = 38
    29 = 619
    30 = 556
}

# Sheet "演出" (Performances)
Set-FValues "演出" @{
    4  = 330
    10 = 1
}

# Sheet "本地生活" (Local Life)
Set-FValues "本地生活" @{
    2 = 466
}

# Sheet "全部类型" (All Types)
Set-FValues "全部类型" @{
    2  = 466
    3  = 7854
    4  = 107
    7  = 13672
    9  = 593
    10 = 648
    11 = 450
    13 = 140
    14 = 441
    15 = 330
    18 = 781
    19 = 44
    20 = 85
    21 = 343
    23 = 20
    26 = 1
    27 = 298
    28 = 141
    29 = 398
    30 = 192
    31 = 1096
    32 = 80
    33 = 641
    34 = 2233
    35 = 754
    36 = 55
    37 = 556
    38 = 38
    40 = 619
    41 = 556
}

Write-Host "Updated want-to-go counts across all sheets."
